$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New values for column J (and mirrored into column N) for rows 2-198
$newValues = @{
    2 = 0.40765130999999999
    3 = 0.38815102000000001
    4 = 0.39850034000000001
    5 = 0.39908130000000003
    6 = 0.39875138999999998
    7 = 0.40786183999999998
    8 = 0.40655913999999999
    9 = 0.39722935999999998
    10 = 0.38965286999999998
    11 = 0.41123833999999998
    12 = 0.38802420999999998
    13 = 0.38574633000000003
    14 = 0.40726688999999999
    15 = 0.39819903000000001
    16 = 0.40277201000000001
    17 = 0.39521838999999997
    18 = 0.40605860999999999
    19 = 0.41318093
    20 = 0.41071361000000001
    21 = 0.40391786000000002
    22 = 0.40696747999999999
    23 = 0.41054139000000001
    24 = 0.41485794999999998
    25 = 0.3999857
    26 = 0.38624484999999997
    27 = 0.39714649000000002
    28 = 0.41408581
    29 = 0.39850966999999998
    30 = 0.39857545
    31 = 0.40340408
    32 = 0.39620242999999999
    33 = 0.39885042999999998
    34 = 0.39678797999999998
    35 = 0.39333125000000002
    36 = 0.39282735000000002
    37 = 0.40702898999999998
    38 = 0.38974545999999999
    39 = 0.39763548999999998
    40 = 0.39822614000000001
    41 = 0.39928152
    42 = 0.40000143999999999
    43 = 0.40171832000000002
    44 = 0.39032641000000001
    45 = 0.39612644000000002
    46 = 0.39595474000000003
    47 = 0.39015212999999999
    48 = 0.40106573000000001
    49 = 0.41314554999999997
    50 = 0.38514469000000001
    51 = 0.39523127000000002
    52 = 0.40022142999999999
    53 = 0.39953915000000001
    54 = 0.39598502000000002
    55 = 0.38250229000000002
    56 = 0.38989900999999999
    57 = 0.40306863999999998
    58 = 0.38410380999999999
    59 = 0.40845673999999998
    60 = 0.42072229999999999
    61 = 0.40708749999999999
    62 = 0.39839174999999999
    63 = 0.39729238
    64 = 0.3974472
    65 = 0.40177480999999998
    66 = 0.39224227
    67 = 0.39697051999999999
    68 = 0.40511434000000002
    69 = 0.39852345
    70 = 0.39329152000000001
    71 = 0.39315464
    72 = 0.38942778
    73 = 0.40552802999999998
    74 = 0.41921638999999999
    75 = 0.39200551
    76 = 0.39940606000000001
    77 = 0.39731031
    78 = 0.41090792999999998
    79 = 0.40271195999999998
    80 = 0.41466752000000001
    81 = 0.38142630999999999
    82 = 0.39071544000000002
    83 = 0.41728517999999998
    84 = 0.40845110000000001
    85 = 0.40975767000000002
    86 = 0.39640482999999999
    87 = 0.40837014999999999
    88 = 0.42720363
    89 = 0.40193335000000002
    90 = 0.40556488000000002
    91 = 0.41043698000000001
    92 = 0.38776632
    93 = 0.39664603999999998
    94 = 0.39839927000000003
    95 = 0.38299074999999999
    96 = 0.39657135999999998
    97 = 0.39185941000000002
    98 = 0.39674903
    99 = 0.39385144999999999
    100 = 0.40145455000000002
    101 = 0.39903569999999999
    102 = 0.38405294000000001
    103 = 0.39665900999999998
    104 = 0.39689906000000003
    105 = 0.40204968000000002
    106 = 0.38707979999999997
    107 = 0.40843244000000001
    108 = 0.40238236999999999
    109 = 0.40427879999999999
    110 = 0.39849970000000001
    111 = 0.40175185000000002
    112 = 0.39225104999999999
    113 = 0.38997256000000002
    114 = 0.41873817000000002
    115 = 0.39164786000000001
    116 = 0.38962744999999999
    117 = 0.40879483
    118 = 0.41027936999999998
    119 = 0.40914274
    120 = 0.39995776
    121 = 0.39027293000000002
    122 = 0.39717622000000002
    123 = 0.38492029999999999
    124 = 0.37997963000000001
    125 = 0.39755531
    126 = 0.40103195000000003
    127 = 0.41117302999999999
    128 = 0.39881486999999999
    129 = 0.38573445000000001
    130 = 0.39952977000000001
    131 = 0.39050615
    132 = 0.40603980000000001
    133 = 0.41069240000000001
    134 = 0.41014721999999998
    135 = 0.40131681000000002
    136 = 0.39644009000000002
    137 = 0.39077748000000001
    138 = 0.37051541999999998
    139 = 0.40917585000000001
    140 = 0.40710622000000002
    141 = 0.41094765
    142 = 0.41267682
    143 = 0.39157511
    144 = 0.39684469999999999
    145 = 0.39592887999999998
    146 = 0.40949746999999997
    147 = 0.40577159000000002
    148 = 0.39173793000000001
    149 = 0.40593194999999999
    150 = 0.40830284
    151 = 0.40912929999999997
    152 = 0.41120641000000002
    153 = 0.40136946000000001
    154 = 0.39742145000000001
    155 = 0.39739022000000002
    156 = 0.40281072000000001
    157 = 0.41190452999999999
    158 = 0.40258830000000001
    159 = 0.38875176
    160 = 0.41077818999999999
    161 = 0.39338641000000002
    162 = 0.40717854999999997
    163 = 0.40072798999999998
    164 = 0.39668695999999998
    165 = 0.39574858000000002
    166 = 0.40390969999999998
    167 = 0.38821724000000002
    168 = 0.39047252999999998
    169 = 0.39709527
    170 = 0.40926591000000001
    171 = 0.38026198999999999
    172 = 0.40406376999999999
    173 = 0.37007909
    174 = 0.41809484000000002
    175 = 0.38587058000000002
    176 = 0.41432885000000003
    177 = 0.39164547999999999
    178 = 0.41529575000000002
    179 = 0.41329634999999998
    180 = 0.38975431999999999
    181 = 0.39643349
    182 = 0.40926593
    183 = 0.40057605000000002
    184 = 0.39796765000000001
    185 = 0.39831357000000001
    186 = 0.39734967999999998
    187 = 0.39363894999999999
    188 = 0.39177255999999999
    189 = 0.39396249
    190 = 0.39909948000000001
    191 = 0.39763999999999999
    192 = 0.39733902999999998
    193 = 0.41900340000000003
    194 = 0.38408063999999997
    195 = 0.42025506000000001
    196 = 0.41968202999999998
    197 = 0.39698240000000001
    198 = 0.39506717000000002
}

foreach ($row in $newValues.Keys) {
    $val = $newValues[$row]
    $ws.Cells.Item($row, 10).Value = $val   # column J
    $ws.Cells.Item($row, 14).Value = $val   # column N
}

# Update the selected cell shown in the sheet view
$ws.Range("R5").Select()